$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 17 (shifts everything below it down by one), then fill it
# with a duplicate of row 16's contents/formatting - a new "Periodo Mora"
# detail line for the same worker (period 2509, in addition to existing 2508).
$ws.Rows("17:17").Insert()

# Copy row 16 formatting onto the newly inserted (blank) row 17
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# Copy row 16 values onto row 17 too
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4163)

# The new row covers the next "Periodo Mora" (2509) instead of 2508
$ws.Range("E17").Value = "2509"

# Update the total "Valor Mora" - now covers 2 periods (56940 + 56940)
$ws.Range("E11").Value = 113880

# Update "Cant. Periodos" - now 2 periods instead of 1
$ws.Range("F13").Value = 2

$excel.CutCopyMode = 0
